# Elton se vc viu isso entao funcionou kkk
#
# Insert a new "Title and Content" slide right after the title slide
# (position 2), pushing every other slide down by one, and fill it in
# with the new test message for Elton.

$p = $ppt.ActivePresentation

# The new slide uses the same "Title and Content" layout as the other
# content slides in this deck (CustomLayouts index 2).
$layout = $p.SlideMaster.CustomLayouts.Item(2)
$newSlide = $p.Slides.AddSlide(2, $layout)

# Title placeholder: "Elton, s vc viu esse slide então funcionou",
# typed as three runs (the middle word "vc" flagged by the spell
# checker), all in Portuguese (Brazil).
$titlePh = $newSlide.Shapes.Placeholders.Item(1)
$titleRange = $titlePh.TextFrame.TextRange
$titleRange.Text = "Elton, s "
$titleRange.LanguageID = "pt-BR"

$run2 = $titleRange.InsertAfter("vc")
$run2.LanguageID = "pt-BR"

$run3 = $titleRange.InsertAfter(" viu esse slide então funcionou")
$run3.LanguageID = "pt-BR"

# Body placeholder: "123"
$bodyPh = $newSlide.Shapes.Placeholders.Item(2)
$bodyRange = $bodyPh.TextFrame.TextRange
$bodyRange.Text = "123"
$bodyRange.LanguageID = "pt-BR"
